$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2024-07-21 Sunday" "2024-07-22 Monday"

Replace-Text "71×25=" "28×55="
Replace-Text "75×38=" "43×85="
Replace-Text "40×39=" "32×41="
Replace-Text "77×63=" "61×88="
Replace-Text "89×37=" "22×80="

Replace-Text "36×22=" "42×41="
Replace-Text "74×71=" "84×29="
Replace-Text "52×91=" "94×44="
Replace-Text "48×22=" "18×83="
Replace-Text "59×73=" "14×23="

Replace-Text "44×82=" "63×84="
Replace-Text "57×85=" "51×45="
Replace-Text "69×39=" "69×35="
Replace-Text "26×68=" "40×52="
Replace-Text "12×78=" "21×20="

Replace-Text "29×89=" "28×99="
Replace-Text "32×30=" "15×89="
Replace-Text "90×37=" "78×97="
Replace-Text "36×84=" "71×68="
Replace-Text "63×94=" "80×94="

Replace-Text "19×18=" "21×90="
Replace-Text "38×69=" "23×32="
Replace-Text "71×98=" "43×75="
Replace-Text "52×76=" "93×49="
Replace-Text "36×85=" "63×31="
